# Fixed test scripts in IAM and Notification Module.
# - Adds a new "IAM029" worksheet (password max-length validation test case)
#   positioned right after "IAM019" and before "Test Case Steps".
# - Adds a corresponding summary row to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add a summary row (row 30) to the "Test Cases" sheet for the new test case
# ---------------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("Test Cases")

# Copy formatting from the row above (row 29) down into the new row 30 so the
# new row picks up the existing borders / wrap styling used throughout the sheet.
$testCases.Range("A29:E29").Copy()
$testCases.Range("A30:E30").PasteSpecial(-4122)

$testCases.Range("A30").Value = "IAM029"
$testCases.Range("B30").Value = "OPQA-2906"
$testCases.Range("C30").Value = "Verify that to validate PASSWORD field in new Neon user registration page with maximum length."
$testCases.Range("D30").Value = "Y"
$testCases.Range("E30").Value = "PASS"

$testCases.Rows.Item(30).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 2) Insert the new "IAM029" worksheet right before "Test Case Steps"
# ---------------------------------------------------------------------------
$testCaseSteps = $wb.Worksheets.Item("Test Case Steps")
$iam029 = $wb.Worksheets.Add($testCaseSteps)
$iam029.Name = "IAM029"

$iam029.Range("A1").Value = "CHARACTER LENGTH"
$iam029.Range("B1").Value = "VALIDITY"
$iam029.Range("C1").Value = "Runmode"
$iam029.Range("D1").Value = "PASS"

$iam029.Range("A2").Value = 91
$iam029.Range("B2").Value = "YES"
$iam029.Range("C2").Value = "Y"
$iam029.Range("D2").Value = "SKIP"

$iam029.Range("A3").Value = 92
$iam029.Range("B3").Value = "YES"
$iam029.Range("C3").Value = "Y"
$iam029.Range("D3").Value = "SKIP"

$iam029.Range("A4").Value = 93
$iam029.Range("B4").Value = "NO"
$iam029.Range("C4").Value = "Y"
$iam029.Range("D4").Value = "PASS"

$iam029.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Fix up selection on the IAM005 sheet (whole-range select, no active cell)
# ---------------------------------------------------------------------------
$iam005 = $wb.Worksheets.Item("IAM005")
$iam005.Activate() | Out-Null
$iam005.Range("A1:D4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Restore the expected active sheet/selection ("Test Cases"!C30)
# ---------------------------------------------------------------------------
$testCases.Activate() | Out-Null
$testCases.Range("C30").Select() | Out-Null
